# Apply the crypto-tracker snapshot refresh captured in the commit:
# prices/volumes updated, and the OKB/Cardano (rows 8-9) and
# Algorand/InternetComputer (rows 37-38) pairs swapped position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "23.806.38"
$ws.Range("E2").Value = "  +2.15%  "

# Row 3
$ws.Range("D3").Value = "1.654.76"
$ws.Range("E3").Value = "  +2.01%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9999"
$ws.Range("E5").Value = "  -0.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.38"
$ws.Range("E6").Value = "  +0.65%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3823"
$ws.Range("E7").Value = "  +1.99%  "

# Row 8
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.32"
$ws.Range("E8").Value = "  -0.18%  "

# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3610"
$ws.Range("E9").Value = "  -0.14%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.248"
$ws.Range("E10").Value = "  +2.38%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08229"
$ws.Range("E11").Value = "  +1.13%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.0000"
$ws.Range("E12").Value = "  -0.18%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.71"
$ws.Range("E13").Value = "  +2.03%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.540"
$ws.Range("E14").Value = "  +1.32%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.423"
$ws.Range("E15").Value = "  +2.37%  "

# Row 16
$ws.Range("E16").Value = "  +0.36%  "

# Row 17
$ws.Range("D17").Value = "1.646.22"
$ws.Range("E17").Value = "  +1.67%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.83"
$ws.Range("E18").Value = "  +4.22%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06979"
$ws.Range("E19").Value = "  +0.76%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.781"
$ws.Range("E20").Value = "  +4.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.80"
$ws.Range("E21").Value = "  +1.69%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.13%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.72"
$ws.Range("E23").Value = "  +2.03%  "

# Row 24
$ws.Range("D24").Value = "23.807.18"
$ws.Range("E24").Value = "  +2.13%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.546"
$ws.Range("E25").Value = "  +3.14%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.102"
$ws.Range("E26").Value = "  +0.76%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.34"
$ws.Range("E27").Value = "  +1.22%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.05"
$ws.Range("E28").Value = "  +0.34%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.242"
$ws.Range("E29").Value = "  -0.48%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.81"
$ws.Range("E30").Value = "  +1.72%  "

# Row 31
$ws.Range("D31").Value = "1.831.74"
$ws.Range("E31").Value = "  +1.89%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.917"
$ws.Range("E32").Value = "  +3.02%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.083"
$ws.Range("E33").Value = "  +1.87%  "

# Row 34
$ws.Range("E34").Value = "  +6.79%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.116"
$ws.Range("E35").Value = "  -2.38%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02852"
$ws.Range("E36").Value = "  +3.96%  "

# Row 37
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.166"
$ws.Range("E37").Value = "  +3.26%  "

# Row 38
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2523"
$ws.Range("E38").Value = "  +1.94%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08843"
$ws.Range("E39").Value = "  +0.81%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07067"
$ws.Range("E40").Value = "  -0.20%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.96"
$ws.Range("E41").Value = "  +8.14%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7091"
$ws.Range("E42").Value = "  +1.93%  "

# Row 43
$ws.Range("E43").Value = "  +0.66%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.92"
$ws.Range("E44").Value = "  -0.27%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6560"
$ws.Range("E45").Value = "  +1.81%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.338"
$ws.Range("E46").Value = "  +3.57%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9994"
$ws.Range("E47").Value = "  -0.13%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.966"
$ws.Range("E48").Value = "  +0.35%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07989"
$ws.Range("E49").Value = "  +0.47%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "129.02"
$ws.Range("E50").Value = "  +2.69%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.197"
$ws.Range("E51").Value = "  +1.35%  "
